$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 87.44388944297727
$ws.Range("C2").Value = 95.3285435439691
$ws.Range("D2").Value = 98.88416789810593
$ws.Range("E2").Value = 98.71753349656692
$ws.Range("F2").Value = 98.35765789288109
$ws.Range("G2").Value = 97.56487881001668
$ws.Range("H2").Value = 96.00027755246491
$ws.Range("B3").Value = 81.68283213612776
$ws.Range("C3").Value = 97.88064996277492
$ws.Range("D3").Value = 99.24955173732829
$ws.Range("E3").Value = 98.64166949886794
$ws.Range("F3").Value = 98.32777591766109
$ws.Range("G3").Value = 97.68121314531946
$ws.Range("H3").Value = 96.34329352983748
$ws.Range("B4").Value = 86.28599674047305
$ws.Range("C4").Value = 95.51137830436474
$ws.Range("D4").Value = 99.0021226259057
$ws.Range("E4").Value = 98.49935476787958
$ws.Range("F4").Value = 98.35125543318377
$ws.Range("G4").Value = 97.5526582025585
$ws.Range("H4").Value = 96.29691818992248
$ws.Range("B5").Value = 84.88664423881967
$ws.Range("C5").Value = 94.4097266602371
$ws.Range("D5").Value = 98.86391237055639
$ws.Range("E5").Value = 98.62439514554762
$ws.Range("F5").Value = 98.266351773633
$ws.Range("G5").Value = 97.39828863319117
$ws.Range("H5").Value = 96.24553532116798
$ws.Range("B6").Value = 85.47579476776683
$ws.Range("C6").Value = 95.43891241958104
$ws.Range("D6").Value = 99.15377616893278
$ws.Range("E6").Value = 98.35709340975923
$ws.Range("F6").Value = 98.26008561321711
$ws.Range("G6").Value = 97.6574570268988
$ws.Range("H6").Value = 96.17813602754524
